# "update site name on google"
# The sheet tracks a cost-per-video table; this adds a new pricing row (row 7)
# for a 146-minute video, following the same formula pattern used by the
# existing rows 3-6, and updates the saved selection/view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Add the new data row (row 7), mirroring rows 3-6 ---
$ws.Range("A7").Value = 146

$ws.Range("B7").Formula = '=0.006*A7'
$ws.Range("C7").Formula = '=B7*1.5'
$ws.Range("D7").Formula = '=C7*100'

$ws.Range("E7").Value = 54000
$ws.Range("F7").Formula = '=15*800'

$ws.Range("G7").Formula = '=(E7*0.00025+F7*0.00125)/1000'
$ws.Range("H7").Formula = '=$K$1*G7'
$ws.Range("I7").Formula = '=H7*100'

# --- Match formatting used by the rows above (thin left/right borders +
#     number formats), same as column B / D for rows 3-6 ---
$ws.Range("B7").NumberFormat = "0.00"
$ws.Range("B7").Interior.ColorIndex = -4142   # xlColorIndexNone (no fill)
$bBorder = $ws.Range("B7").Borders.Item(7)    # xlEdgeLeft
$bBorder.Color = 0
$bBorder.LineStyle = 1                        # xlContinuous
$bBorder.Weight = 2                           # xlThin

$ws.Range("C7").NumberFormat = "0.00"

$ws.Range("D7").NumberFormat = "0.0"
$ws.Range("D7").Interior.ColorIndex = -4142
$dBorder = $ws.Range("D7").Borders.Item(10)   # xlEdgeRight
$dBorder.Color = 0
$dBorder.LineStyle = 1
$dBorder.Weight = 2

# --- Update the saved selection / scroll position for the sheet ---
$ws.Range("O9").Select()

$wb.Save()
